$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the plan-content text in column B to reflect the author's edits.
$ws.Range("B3").Value = "1、熟悉需求规格说明书；2、确定pc端用例（参与者为<普通用户>，功能模块为<用户群组管理>）"
$ws.Range("B4").Value = "1、熟悉需求规格说明书；2、确定pc端用例（参与者为<普通用户>，功能模块为<用户群组管理>）"
$ws.Range("B5").Value = "1、熟悉APP手机端UI设计；2、确定用例（参与者为<普通用户>，功能模块为<用户管理>）最后汇总"
$ws.Range("B6").Value = "1、熟悉Android项目-Java项目整合开发-天下纵横-系统需求分析-吴绍根；2、确定用例（参与者为<普通用户>，功能模块为<用户管理>）最后汇总"

# Update the active selection on the sheet.
$ws.Range("M13").Select()
